$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105; all existing rows 105-154 shift down to 106-155.
$ws.Rows.Item(105).EntireRow.Insert()

# Populate the newly inserted row 105 with the new weekly price record.
$ws.Cells.Item(105, 1).Value  = 11
$ws.Cells.Item(105, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(105, 3).Value  = "Bíobío"
$ws.Cells.Item(105, 4).Value  = 44636
$ws.Cells.Item(105, 5).Value  = 8
$ws.Cells.Item(105, 6).Value  = "Fruta"
$ws.Cells.Item(105, 7).Value  = 100108
$ws.Cells.Item(105, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(105, 9).Value  = 100108005
$ws.Cells.Item(105, 10).Value = "Piña"
$ws.Cells.Item(105, 11).Value = "Caramelo"
$ws.Cells.Item(105, 12).Value = "Segunda"
$ws.Cells.Item(105, 13).Value = 220
$ws.Cells.Item(105, 14).Value = 17000
$ws.Cells.Item(105, 15).Value = 18000
$ws.Cells.Item(105, 16).Value = 17455
$ws.Cells.Item(105, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(105, 18).Value = "Ecuador"
$ws.Cells.Item(105, 19).Value = 1247
$ws.Cells.Item(105, 20).Value = 14
